# "Generate Report for Archive" — regenerate the localization-status report:
#   * every cell whose text is "Ready for handoff" becomes "In Translation"
#     (status text changed on the Overview sheet + the per-locale sheets)
#   * the Status columns that held the (now shorter) text get narrower,
#     matching the width Excel's own column auto-fit would produce.

$wb = $excel.ActiveWorkbook

$oldText = "Ready for handoff"
$newText = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # Put the literal on the left of -eq: some cells hold real
            # booleans (True/False) and comparing "$bool -eq $string" would
            # coerce the string instead of the other way round.
            if ($oldText -eq $cell.Value2) {
                $cell.Value = $newText
            }
        }
    }
}

# Narrow the columns that used to show "Ready for handoff" so they fit the
# shorter "In Translation" text, same as Excel would do when it reflows the
# report.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E1").ColumnWidth = 12.5
$ws1.Range("F1").ColumnWidth = 12.5

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C1").ColumnWidth = 12.5

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C1").ColumnWidth = 12.5
